# Update the "Förändrad" (Changed) date column (C) for rows 2-20
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C20").Value = 45185
